$d = $word.ActiveDocument

# wdReplace constants used below:
#   wdFindContinue / wdReplaceAll = 2, wdReplaceOne = 1

# --- 1. Title (Heading1) + the bold "title" run near the bottom: both share the same
#        old text, so a single ReplaceAll across the whole document content covers both. ---
$d.Content.Find.Execute(
    "Play Firebird Double 27 Free: Review and Ratings", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Firebird Double 27 for Free", 2) | Out-Null

# --- 2. Insert a brand-new bullet "Simple and intuitive gameplay" right before
#        "Up to 54 ways to win" in the "What we like" list. ---
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Up to 54 ways to win") {
        $newPara = $p.Range.InsertParagraphBefore()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "") {
        $nxt = $p.Next()
        if ($nxt -ne $null) {
            $nxtTxt = $nxt.Range.Text.TrimEnd([char]13, [char]7)
            if ($nxtTxt -eq "Up to 54 ways to win") {
                $p.Range.Text = "Simple and intuitive gameplay"
                break
            }
        }
    }
}

# --- 3. Bullet text swaps in "What we like" ---
$d.Content.Find.Execute(
    "Wild Phoenix symbol can replace all symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Wild symbol substitutes for other symbols", 2) | Out-Null

# --- 4. Drop the "Intuitive gameplay for beginners" bullet entirely ---
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Intuitive gameplay for beginners") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- 5. "Big jackpots available" -> "Decent jackpots and winning opportunities" ---
$d.Content.Find.Execute(
    "Big jackpots available", $true, $false, $false, $false, $false,
    $true, 1, $false, "Decent jackpots and winning opportunities", 2) | Out-Null

# --- 6. "What we don't like" bullet text swaps ---
$d.Content.Find.Execute(
    "Lacks advanced features", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lack of additional special features", 2) | Out-Null

$d.Content.Find.Execute(
    "May not appeal to experienced players", $true, $false, $false, $false, $false,
    $true, 1, $false, "May not appeal to advanced players", 2) | Out-Null

# --- 7. Meta description paragraph at the very end ---
$d.Content.Find.Execute(
    "Learn all about Firebird Double 27 slot game. Play for free and read our review of this traditional fruit-themed slot with big jackpots available.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Firebird Double 27, a simple and intuitive slot game with decent jackpots, for free.", 2) | Out-Null
